$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Exportado:" timestamp shown in D3 ---
$ws.Range("D3").Value = "mié. 07/05/2025 19:09"

# --- Grow the "Datos" table with 4 more daily rows, pushing the totals row down ---
# New data rows: date serial, Cantidad de pedidos, Total $, Cantidad de productos
$newRows = @(
  @(45780, 19, 418500, 30),
  @(45781, 12, 201200, 23),
  @(45783, 14, 275800, 37),
  @(45784, 3, 55500, 4)
)

$firstNewRow = 38

for ($i = 0; $i -lt $newRows.Count; $i++) {
  $r = $firstNewRow + $i
  # Insert a fresh row right above the (old) totals row; this shifts the
  # totals row (and its formulas/styles) down by one each time.
  $ws.Rows.Item($r).Insert()

  $row = $newRows[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
}

$newTotalsRow = $firstNewRow + $newRows.Count

# --- Resize the table/autofilter to cover the new range ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A5:D" + $newTotalsRow))

Write-Host ("Table range: " + $lo.Range.Address())
Write-Host "Done"
